# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" column for the rows that are
# currently pending handoff (status "Ready for handoff" / "Handback
# transform failed") so that they reflect the freshly generated handoff
# timestamp for each locale.

$wb = $excel.ActiveWorkbook

$zhRows = @(7, 10, 11, 12, 13, 14, 15, 16)
$deRows = @(7, 10, 11, 12, 13, 14, 15, 16)

$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $zhRows) {
    $wsZh.Range("D$r").Value = "2016-03-09 06:35:32"
}

$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $deRows) {
    $wsDe.Range("D$r").Value = "2016-03-09 06:35:36"
}
